$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.321.93'
$ws.Range("E2").Value = '  +0.44%  '

$ws.Range("D3").Value = '3.153.42'
$ws.Range("E3").Value = '  -1.22%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.88%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.63%  '

$ws.Range("E9").Value = '  -3.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.382'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("D12").Value = '3.705.82'
$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("E13").Value = '  -1.02%  '

$ws.Range("D14").Value = '64.407.87'
$ws.Range("E14").Value = '  +0.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.88%  '

$ws.Range("D16").Value = '3.157.22'
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("E17").Value = '  -2.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '404.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.60%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.37'
$ws.Range("D23").Style = "Normal"

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.482'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.81%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.195'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.36%  '

$ws.Range("E26").Value = '  -6.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("E35").Value = '  -2.84%  '

$ws.Range("D36").Value = '2.679.04'
$ws.Range("E36").Value = '  -2.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.696'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0617'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.91%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.95%  '

$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '290.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.11%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0257'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.11%  '

$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("E47").Value = '  -1.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.878'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.46%  '
